$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '20.824.95'
$ws.Range("E2").Value = '  +2.28%  '
$ws.Range("D3").Value = '1.517.08'
$ws.Range("E3").Value = '  +4.15%  '
$ws.Range("D4").Value = '''1.003'
$ws.Range("E4").Value = '  -0.23%  '
$ws.Range("D5").Value = '''0.9592'
$ws.Range("E5").Value = '  +2.05%  '
$ws.Range("D6").Value = '''279.79'
$ws.Range("E6").Value = '  +2.03%  '
$ws.Range("D7").Value = '''0.3585'
$ws.Range("E7").Value = '  -1.68%  '
$ws.Range("D8").Value = '''0.3130'
$ws.Range("E8").Value = '  +1.95%  '
$ws.Range("D9").Value = '''1.105'
$ws.Range("E9").Value = '  +7.04%  '
$ws.Range("D10").Value = '''39.86'
$ws.Range("E10").Value = '  -0.18%  '
$ws.Range("D11").Value = '''0.06733'
$ws.Range("E11").Value = '  +2.99%  '
$ws.Range("D12").Value = '''0.9973'
$ws.Range("E12").Value = '  +0.09%  '
$ws.Range("D13").Value = '''18.62'
$ws.Range("E13").Value = '  +4.98%  '
$ws.Range("D14").Value = '''5.588'
$ws.Range("E14").Value = '  +3.89%  '
$ws.Range("D15").Value = '''6.260'
$ws.Range("E15").Value = '  +2.79%  '
$ws.Range("D16").Value = '''0.9632'
$ws.Range("E16").Value = '  +0.68%  '
$ws.Range("D17").Value = '''0.00001030'
$ws.Range("E17").Value = '  +1.01%  '
$ws.Range("D18").Value = '1.506.95'
$ws.Range("E18").Value = '  +3.82%  '
$ws.Range("D19").Value = '''0.06038'
$ws.Range("E19").Value = '  +5.70%  '
$ws.Range("D20").Value = '''70.31'
$ws.Range("E20").Value = '  +1.34%  '
$ws.Range("D21").Value = '''5.602'
$ws.Range("E21").Value = '  +3.43%  '
$ws.Range("D22").Value = '''14.87'
$ws.Range("E22").Value = '  +3.39%  '
$ws.Range("D23").Value = '''11.34'
$ws.Range("E23").Value = '  +4.68%  '
$ws.Range("D24").Value = '''2.309'
$ws.Range("E24").Value = '  +2.82%  '
$ws.Range("D25").Value = '20.844.72'
$ws.Range("E25").Value = '  +2.38%  '
$ws.Range("D26").Value = '''146.42'
$ws.Range("E26").Value = '  +4.19%  '
$ws.Range("D27").Value = '''2.154'
$ws.Range("E27").Value = '  +2.80%  '
$ws.Range("D28").Value = '''17.48'
$ws.Range("E28").Value = '  +2.36%  '
$ws.Range("D29").Value = '1.670.51'
$ws.Range("E29").Value = '  +4.22%  '
$ws.Range("D30").Value = '''115.93'
$ws.Range("D31").Value = '''4.008'
$ws.Range("E31").Value = '  +2.07%  '
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").Value = '''0.8371'
$ws.Range("E32").Value = '  +5.83%  '
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = '''5.086'
$ws.Range("E33").Value = '  +4.98%  '
$ws.Range("D34").Value = '''0.07993'
$ws.Range("E34").Value = '  +3.15%  '
$ws.Range("B35").Value = 'WEMIXTOKEN'
$ws.Range("C35").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D35").Value = '''1.474'
$ws.Range("E35").Value = '  -1.46%  '
$ws.Range("B36").Value = 'TrustWalletToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D36").Value = '''1.210'
$ws.Range("E36").Value = '  +7.35%  '
$ws.Range("B37").Value = 'InternetComputer(DFINITY)'
$ws.Range("C37").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D37").Value = '''4.851'
$ws.Range("E37").Value = '  +3.95%  '
$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").Value = '''0.05793'
$ws.Range("E38").Value = '  +1.95%  '
$ws.Range("D39").Value = '''0.02075'
$ws.Range("E39").Value = '  +2.91%  '
$ws.Range("D40").Value = '''10.54'
$ws.Range("E40").Value = '  +1.86%  '
$ws.Range("D41").Value = '''0.9621'
$ws.Range("E41").Value = '  +1.55%  '
$ws.Range("D42").Value = '''0.1884'
$ws.Range("E42").Value = '  +1.42%  '
$ws.Range("D43").Value = '''7.534'
$ws.Range("E43").Value = '  +2.29%  '
$ws.Range("D44").Value = '''0.5328'
$ws.Range("E44").Value = '  +1.61%  '
$ws.Range("D45").Value = '''3.543'
$ws.Range("E45").Value = '  +1.77%  '
$ws.Range("D46").Value = '''12.24'
$ws.Range("E46").Value = '  +2.87%  '
$ws.Range("D47").Value = '''120.47'
$ws.Range("E47").Value = '  +2.78%  '
$ws.Range("D48").Value = '''0.5329'
$ws.Range("E48").Value = '  +4.01%  '
$ws.Range("D49").Value = '''1.857'
$ws.Range("E49").Value = '  +6.73%  '
$ws.Range("D50").Value = '''0.06504'
$ws.Range("E50").Value = '  +1.46%  '
$ws.Range("D51").Value = '''0.9863'
$ws.Range("E51").Value = '  +0.58%  '
